$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values (Price/Volume columns) must be written with a leading
# apostrophe so Excel keeps them as literal text (matching the source inlineStr
# cells) instead of silently parsing them into Number/Percentage values. The
# Style reset afterwards clears the "quote prefix" text style Excel applies so
# the cell keeps its original (unstyled) formatting.

$ws.Range("D2").Value = "'325.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.18%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.05%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.681"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'7.21%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08033"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.84%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.028"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.46%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.486"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'8.620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.41%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-1.48%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9227"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.19%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1244"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-8.57%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1962"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'8.723"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'20.73%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09209"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.83%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.03562"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.19%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'9.34%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.001300"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.41%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006091"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-5.93%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.350"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.43%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3481"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.85%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1370"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.85%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2503"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.28%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04376"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.34%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001261"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.19%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'6.17%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'2.46%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.57%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05332"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.42%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007480"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.30%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009906"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'9.30%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1406"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.60%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002116"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.12%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.62%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006678"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.33%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002279"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-5.09%"
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-11.02%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.10%"
$ws.Range("E51").Style = "Normal"
